$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.369.98"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").Value = "2.081.85"
$ws.Range("E3").Value = "  +2.32%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0834"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").Value = "2.389.48"
$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.94%  "

$ws.Range("E15").Value = "  +1.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.56%  "

$ws.Range("D17").Value = "2.081.88"
$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").Value = "38.337.69"
$ws.Range("E18").Value = "  +1.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.13%  "

$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +1.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.52%  "

$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.91%  "

$ws.Range("E34").Value = "  +2.86%  "

$ws.Range("E35").Value = "  +0.43%  "

$ws.Range("E36").Value = "  +2.28%  "

$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.28%  "

$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("E40").Value = "  +2.20%  "

$ws.Range("D41").Value = "1.538.98"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.30%  "

$ws.Range("E44").Value = "  +1.65%  "

$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.14%  "

$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").Value = "2.276.75"
$ws.Range("E51").Value = "  +2.32%  "

